$p = $ppt.ActivePresentation

$deg = [char]0x00B0   # ° degree sign
$sup2 = [char]0x00B2  # ² superscript two

# --- Slide 1 ---
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction to In4SnS8 Nanosheets"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Overview of In4SnS8 nanosheets`r- Applications in environmental remediation, solar energy conversion, and advanced nanodevices`r- Importance in nanotechnology and material science"

# --- Slide 2 ---
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Synthesis of In4SnS8 Nanosheets"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Synthesized via a thermal decomposition method`r- Utilizes Sn(DDTC)4 and In(DDTC)3 in oleylamine solvent`r- Heating process: 120${deg}C under vacuum, followed by 240${deg}C under N2 atmosphere"

# --- Slide 3 ---
$s = $p.Slides.Item(3)
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Morphology examined using TEM (Tecnai G2 F30 S-Twin)`r- XRD patterns for structural analysis (Bruker D8 Advance)`r- XPS for elemental composition (PHI 5000 Versaprobe)"

# --- Slide 4 ---
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Properties of In4SnS8 Nanosheets"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Ultrathin nature with an average thickness of 3.8 nm`r- Comprises five atomically thick layers`r- Large specific surface area of 40.34 m${sup2}/g"

# --- Slide 5 ---
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Photocatalytic Efficiency"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Fast adsorption-visible-light photocatalysis dual function`r- Effective for various organic dyes removal`r- Potential for solar energy conversion and environmental remediation"

# --- Slide 6 ---
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Comparative Analysis"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Higher efficiency compared to other structures like In4SnS8 microspheres`r- Specific surface area significantly larger than flower-like microspheres (24.7 m${sup2}/g)"

# --- Slide 7 ---
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion and Future Directions"
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = " "
$body.Text = "- Promising advances in 2D nanomaterials`r- Potential applications in environmental and energy sectors`r- Continued research to optimize synthesis and application methods"
